$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.0292345
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.4428765120700495
$ws.Range("J2").Value = 0.346386487911515
$ws.Range("O2").Value = 0.4086672402490986
$ws.Range("P2").Value = 0.5089958879585649
$ws.Range("Q2").Value = 0.0011664273155
$ws.Range("R2").Value = 0.006998563893
$ws.Range("S2").Value = 0.1809891219588137
$ws.Range("T2").Value = 0.1763092979913703

# Row 3
$ws.Range("G3").Value = 0.0292345
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.4428765120700495
$ws.Range("J3").Value = 0.346386487911515
$ws.Range("M3").Value = 0.057733
$ws.Range("N3").Value = 0.115466
$ws.Range("O3").Value = 0.5913327597509014
$ws.Range("P3").Value = 0.4910041120414351
$ws.Range("Q3").Value = 0.0016877953885
$ws.Range("R3").Value = 0.006751181554
$ws.Range("S3").Value = 0.2618873901112357
$ws.Range("T3").Value = 0.1700771899201447

# Row 4
$ws.Range("I4").Value = 0.5571234879299505
$ws.Range("J4").Value = 0.6536135120884849
$ws.Range("O4").Value = 0.4086672402490986
$ws.Range("P4").Value = 0.5089958879585649
$ws.Range("S4").Value = 0.2276781182902849
$ws.Range("T4").Value = 0.3326865899671946

# Row 5
$ws.Range("I5").Value = 0.5571234879299505
$ws.Range("J5").Value = 0.6536135120884849
$ws.Range("M5").Value = 0.057733
$ws.Range("N5").Value = 0.115466
$ws.Range("O5").Value = 0.5913327597509014
$ws.Range("P5").Value = 0.4910041120414351
$ws.Range("Q5").Value = 0.002123188808
$ws.Range("R5").Value = 0.012739132848
$ws.Range("S5").Value = 0.3294453696396656
$ws.Range("T5").Value = 0.3209269221212903
